$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Test Cases")

# Increment the "Line Number" values (column C) in rows 5 through 15 by 1,
# reflecting updated line numbers after stabilizing generate_response logic.
for ($row = 5; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
